$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.339.49"
$ws.Range("E2").Value = "  -0.07%  "
Set-TextValue $ws.Range("D3") "1.933.83"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "0.7453"
$ws.Range("E5").Value = "  +2.58%  "
Set-TextValue $ws.Range("D6") "248.92"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "28.09"
$ws.Range("E8").Value = "  -1.37%  "
Set-TextValue $ws.Range("D9") "0.3200"
$ws.Range("E9").Value = "  -4.15%  "
Set-TextValue $ws.Range("D10") "0.07110"
$ws.Range("E10").Value = "  -2.00%  "
Set-TextValue $ws.Range("D11") "0.7882"
$ws.Range("E11").Value = "  -2.89%  "
Set-TextValue $ws.Range("D12") "0.08002"
$ws.Range("E12").Value = "  -1.20%  "
Set-TextValue $ws.Range("D13") "1.933.24"
$ws.Range("E13").Value = "  -0.33%  "
Set-TextValue $ws.Range("D14") "5.379"
$ws.Range("E14").Value = "  -1.87%  "
Set-TextValue $ws.Range("D15") "94.49"
$ws.Range("E15").Value = "  -0.20%  "
Set-TextValue $ws.Range("D16") "14.61"
$ws.Range("E16").Value = "  -2.81%  "
Set-TextValue $ws.Range("D17") "30.344.30"
$ws.Range("E17").Value = "  -0.06%  "
Set-TextValue $ws.Range("D18") "252.54"
$ws.Range("E18").Value = "  +0.56%  "
Set-TextValue $ws.Range("D19") "0.000008032"
$ws.Range("E19").Value = "  -2.51%  "
Set-TextValue $ws.Range("D20") "5.783"
$ws.Range("E20").Value = "  -2.32%  "
Set-TextValue $ws.Range("D21") "2.188.24"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.06%  "
Set-TextValue $ws.Range("D24") "6.823"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -2.18%  "
Set-TextValue $ws.Range("D26") "164.50"
$ws.Range("E26").Value = "  +0.85%  "
Set-TextValue $ws.Range("D27") "2.325"
$ws.Range("E27").Value = "  -2.98%  "
Set-TextValue $ws.Range("D28") "19.11"
$ws.Range("E28").Value = "  -1.16%  "
Set-TextValue $ws.Range("D29") "0.1316"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  +0.38%  "
Set-TextValue $ws.Range("D31") "1.531"
$ws.Range("E31").Value = "  -2.69%  "
Set-TextValue $ws.Range("D32") "4.435"
$ws.Range("E32").Value = "  -0.14%  "
Set-TextValue $ws.Range("D33") "4.149"
$ws.Range("E33").Value = "  -1.48%  "
Set-TextValue $ws.Range("D34") "0.05137"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  -0.94%  "
Set-TextValue $ws.Range("D36") "0.7487"
$ws.Range("E36").Value = "  -0.42%  "
Set-TextValue $ws.Range("D37") "2.770"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  -0.84%  "
Set-TextValue $ws.Range("D39") "2.813"
$ws.Range("E39").Value = "  -0.94%  "
Set-TextValue $ws.Range("D40") "78.01"
$ws.Range("E40").Value = "  -3.70%  "
Set-TextValue $ws.Range("D41") "6.411"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("E43").Value = "  -2.75%  "
Set-TextValue $ws.Range("D44") "0.8436"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  +0.04%  "
Set-TextValue $ws.Range("D46") "102.32"
$ws.Range("E46").Value = "  -0.05%  "
Set-TextValue $ws.Range("D47") "9.827"
$ws.Range("E47").Value = "  +0.00%  "
Set-TextValue $ws.Range("D48") "7.529"
$ws.Range("E48").Value = "  +0.86%  "

# Row 49/50 swap: Maker <-> Elrond
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D49") "37.47"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D50") "984.97"
$ws.Range("E50").Value = "  +11.88%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.1193"
$ws.Range("E51").Value = "  +4.43%  "

Write-Output "applied"